$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "SELECT `n    LOCATIONID, `n    COUNT(*) AS NumberOfTurtlesMeasured`nFROM `n    TBLFIELDDATATURTLEMEASUREMENTS`nGROUP BY `n    LOCATIONID"
$ws.Range("I3").Value = "SELECT `n    LOCATIONID, `n    COUNT(*) AS NumberOfTurtlesMeasured`nFROM `n    TBLFIELDDATATURTLEMEASUREMENTS`nGROUP BY `n    LOCATIONID"
$ws.Range("H12").Value = "SELECT `n    TMT.RECORDID AS TurtleMeasurementRecordID, `n    TMT.COMMENTS AS MeasurementComments`nFROM `n    TBLFIELDDATATURTLEMEASUREMENTS TMT`nJOIN `n    TBLFIELDDATATURTLETRAPSURVEYS TTS ON TMT.EVENTID = TTS.EVENTID AND TMT.LOCATIONID = TTS.LOCATIONID`nWHERE `n    TMT.SPECIES_CODE = TTS.SPECIES_CODE `n    AND TMT.NOTCH_CODE = TTS.NOTCH_CODE `n    AND TMT.SEX <> TTS.SEX"
$ws.Range("I12").Value = "SELECT `n    TMT.RECORDID AS TurtleMeasurementRecordID, `n    TMT.COMMENTS AS MeasurementComments`nFROM `n    TBLFIELDDATATURTLEMEASUREMENTS TMT`nJOIN `n    TBLFIELDDATATURTLETRAPSURVEYS TTS ON TMT.EVENTID = TTS.EVENTID AND TMT.LOCATIONID = TTS.LOCATIONID`nWHERE `n    TMT.SPECIES_CODE = TTS.SPECIES_CODE `n    AND TMT.NOTCH_CODE = TTS.NOTCH_CODE `n    AND TMT.SEX <> TTS.SEX"
$ws.Range("H17").Value = "SELECT o.[AGENCY/TITLE] AS Agency, COUNT(t.RECORDID) AS TurtleMeasurementCount`nFROM TBLFIELDDATATURTLEMEASUREMENTS t`nJOIN TBLEVENTDATAHERPS e ON t.EVENTID = e.EVENTID`nJOIN OBSERVER_LU o ON e.OBSINITS = o.OBSINITS`nGROUP BY o.[AGENCY/TITLE]"
$ws.Range("I17").Value = "SELECT o.[AGENCY/TITLE] AS Agency, COUNT(t.RECORDID) AS TurtleMeasurementCount`nFROM TBLFIELDDATATURTLEMEASUREMENTS t`nJOIN TBLEVENTDATAHERPS e ON t.EVENTID = e.EVENTID`nJOIN OBSERVER_LU o ON e.OBSINITS = o.OBSINITS`nGROUP BY o.[AGENCY/TITLE]"
$ws.Range("H19").Value = "SELECT DISTINCT O.OBSINITS, O.FIRSTNAME, O.LASTNAME`nFROM OBSERVER_LU O`nJOIN TLINKOBSERVERS LO ON O.OBSINITS = LO.OBSINITS`nJOIN TBLFIELDDATATURTLEMEASUREMENTS TM ON LO.EVENTID = TM.EVENTID"
$ws.Range("I19").Value = "SELECT DISTINCT O.OBSINITS, O.FIRSTNAME, O.LASTNAME`nFROM OBSERVER_LU O`nJOIN TLINKOBSERVERS LO ON O.OBSINITS = LO.OBSINITS`nJOIN TBLFIELDDATATURTLEMEASUREMENTS TM ON LO.EVENTID = TM.EVENTID"
$ws.Range("H22").Value = "SELECT S.DESCRIPTION AS StageDescription, COUNT(MT.RECORDID) AS MinnowCount`nFROM TBLFIELDDATAMINNOWTRAPSURVEYS MT`nJOIN TLUSTAGE S ON MT.STAGE = S.STAGE`nGROUP BY S.DESCRIPTION`nORDER BY MinnowCount DESC"
$ws.Range("I22").Value = "SELECT S.DESCRIPTION AS StageDescription, COUNT(MT.RECORDID) AS MinnowCount`nFROM TBLFIELDDATAMINNOWTRAPSURVEYS MT`nJOIN TLUSTAGE S ON MT.STAGE = S.STAGE`nGROUP BY S.DESCRIPTION`nORDER BY MinnowCount DESC"
$ws.Range("H27").Value = "SELECT B.BEHAVIOR, COUNT(*) AS RecordCount`nFROM TBLFIELDDATASNAKEDATACOLLECTION SDC`nJOIN TLUBEHAVIOR B ON SDC.BEHAVIOR = B.BEHAVIOR`nJOIN TBLFIELDDATACOVERBOARD CB ON SDC.EVENTID = CB.EVENTID AND SDC.SNAKEID = CB.SNAKEID`nWHERE CB.TYPE = 'coverboard'`nGROUP BY B.BEHAVIOR"
$ws.Range("I27").Value = "SELECT B.BEHAVIOR, COUNT(*) AS RecordCount`nFROM TBLFIELDDATASNAKEDATACOLLECTION SDC`nJOIN TLUBEHAVIOR B ON SDC.BEHAVIOR = B.BEHAVIOR`nJOIN TBLFIELDDATACOVERBOARD CB ON SDC.EVENTID = CB.EVENTID AND SDC.SNAKEID = CB.SNAKEID`nWHERE CB.TYPE = 'coverboard'`nGROUP BY B.BEHAVIOR"
$ws.Range("H29").Value = "SELECT `n    LP.POINTID, `n    FDCB.SNAKEID, `n    FDCB.[BOARD_"
$ws.Range("I29").Value = "SELECT `n    LP.POINTID, `n    FDCB.SNAKEID, `n    FDCB.[BOARD_"
$ws.Range("H31").Value = "SELECT `n    WEATHER, `n    COUNT(EVENTID) AS EventCount, `n    AVG(AIRTEMP) AS AvgAirTemp, `n    AVG(WATERTEMP) AS AvgWaterTemp`nFROM `n    TBLEVENTDATAHERPS`nGROUP BY `n    WEATHER"
$ws.Range("I31").Value = "SELECT `n    WEATHER, `n    COUNT(EVENTID) AS EventCount, `n    AVG(AIRTEMP) AS AvgAirTemp, `n    AVG(WATERTEMP) AS AvgWaterTemp`nFROM `n    TBLEVENTDATAHERPS`nGROUP BY `n    WEATHER"
$ws.Range("H32").Value = "SELECT `n    EDH.OBSINITS, `n    AVG(EDH.AIRTEMP) AS AverageAirTemperature`nFROM `n    TBLEVENTDATAHERPS AS EDH`nGROUP BY `n    EDH.OBSINITS"
$ws.Range("I32").Value = "SELECT `n    EDH.OBSINITS, `n    AVG(EDH.AIRTEMP) AS AverageAirTemperature`nFROM `n    TBLEVENTDATAHERPS AS EDH`nGROUP BY `n    EDH.OBSINITS"
$ws.Range("H33").Value = "SELECT `n    EDH.WEATHER, `n    AVG(EDH.AIRTEMP) AS AverageAirTemperature`nFROM `n    TBLEVENTDATAHERPS AS EDH`nGROUP BY `n    EDH.WEATHER"
$ws.Range("I33").Value = "SELECT `n    EDH.WEATHER, `n    AVG(EDH.AIRTEMP) AS AverageAirTemperature`nFROM `n    TBLEVENTDATAHERPS AS EDH`nGROUP BY `n    EDH.WEATHER"
